$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" footer date from 3/23/2023 to
#    4/4/2023 everywhere it appears (Slide Master + every Slide Layout).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($container, $newText) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "4/4/2023"

# Slide Master
Update-DatePlaceholder $p.SlideMaster $newDate

# Every Slide Layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j) $newDate
}

# Notes Master (best effort; harmless if unsupported)
try {
    Update-DatePlaceholder $p.NotesMaster $newDate
} catch {
}

# ---------------------------------------------------------------------------
# 2) Give the vertical "Tab" shapes extra top margin/inset: tIns 182880 EMU
#    (14.4pt) -> 274320 EMU (21.6pt) on every "Tab N" shape inside the
#    "Haiku" group on every slide.
# ---------------------------------------------------------------------------
$newMarginTopPt = 21.6

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $top = $slide.Shapes.Item($shi)
        if ($top.Name -eq "Haiku") {
            for ($gi = 1; $gi -le $top.GroupItems.Count; $gi++) {
                $item = $top.GroupItems.Item($gi)
                if ($item.Name -like "Tab *") {
                    $item.TextFrame.MarginTop = $newMarginTopPt
                }
            }
        }
    }
}
